$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (registro 2): nova data de entrega e quantidades zeradas
$ws.Range("C3").Value = "20/01/2023"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0

# Row 4 (registro 3): concluido, com valores preenchidos
# (formato texto evita que "02/03/2021" vire numero de serie de data)
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "02/03/2021"
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 50
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 25805
$ws.Range("K4").Value = "Concluído"

# Row 5 (registro 4): nome e data atualizados
$ws.Range("B5").Value = "Jean Carlos"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "03/01/2022"
$ws.Range("E5").Value = 1

# Row 6 (registro 5): data/hora atualizadas, valores zerados, status pendente
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "03/01/2022"
$ws.Range("D6").Value = "10:30"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").ClearContents()
$ws.Range("J6").Font.Bold = $false
$ws.Range("K6").Value = "Pendente"
